$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_6.3")

$months = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($row = 5; $row -le 84; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $monthNum = [int]$cell.Value()
    $cell.Value = $months[$monthNum]
}

# Writing into hidden (filtered-out) rows causes the engine to stamp an
# explicit row height on them as a side effect. Restore each touched
# hidden row to its pre-edit height state so only the cell contents
# change, matching the original row layout.
$rowsWithExplicitHeight = @(37, 84)
for ($row = 25; $row -le 84; $row++) {
    if ($rowsWithExplicitHeight -contains $row) {
        $ws.Rows.Item($row).RowHeight = 18
    } else {
        $ws.Rows.Item($row).AutoFit()
    }
}
